$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the "Calculations" section (old row 8 -> new row 10)
$ws.Rows("8:9").Insert()

# Update the input label and value
$ws.Range("A2").Value = "Total Price Before Tax"
$ws.Range("B2").Value = 113.95

# Update "Total After Tax" formula to reference the new total row (B13)
$ws.Range("B7").Formula = "=B13"

# Add the new "Total Amount round upto 2 decimal" row at the bottom
$ws.Range("A13").Value = "Total Amount round upto 2 decimal"
$ws.Range("B13").Formula = "=ROUND(B2+B12,2)"
$ws.Range("B13").Style = $ws.Range("B12").Style

# Adjust column width for column A
$ws.Columns("A").ColumnWidth = 31.5546875

# Update the active selection
$ws.Range("A3").Select()
